# Updates IG sushi configuration
#
# The "Metadata" worksheet (sheet 1) is a Property/Value table generated by
# the FHIR IG publisher ("sushi"). This regeneration:
#   - refreshes the generation Date and the Publisher name
#   - adds a new "Jurisdiction" / "Germany" row right after the existing
#     "Contact" row (which ends up duplicated immediately above it, exactly
#     as the regenerated table does)
#   - every row below shifts down by two rows to make room

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the Date and Publisher values (row 8 / row 9) ---
$ws.Range("B8").Value = "2021-12-22T21:26:07+01:00"
$ws.Range("B9").Value = "Forschungsgruppe Digital Health"

# --- Make room for the duplicated "Contact" row and the new "Jurisdiction" row ---
# Row 10 is "Contact" / "No display for ContactDetail". Insert two blank rows
# right after it (new rows 11 and 12), shifting everything from the old row 11
# onward down by two.
$ws.Rows.Item(11).Resize(2).Insert()

# Copy the formatting (borders/fill/font/alignment) of row 10 onto the two new
# rows so they keep the same "data row" style as every other row.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B12").PasteSpecial(-4122)

# New row 11: duplicate of the "Contact" row.
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "No display for ContactDetail"

# New row 12: the new "Jurisdiction" / "Germany" entry.
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = "Germany"
